$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Add the new "InlineNote" character style (based on Default
#    Paragraph Font, 10pt / sz=20, ui priority 1, quick style).
# ------------------------------------------------------------------
$inlineNote = $d.Styles.Add("InlineNote", 2)
$inlineNote.BaseStyle = $d.Styles("DefaultParagraphFont")
$inlineNote.Priority = 1
$inlineNote.Font.Size = 10
$inlineNote.QuickStyle = $true

# ------------------------------------------------------------------
# 2. Locate the empty "Note"-styled paragraph that sits just before
#    the manual page break (the other empty "Note" paragraph, right
#    before the demo table, must stay untouched).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$target = $null
$fallback = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Style.NameLocal -eq "Note" -and $p.Range.Text.Trim().Length -eq 0) {
        if ($fallback -eq $null) {
            $fallback = $p
        }
        if ($i -lt $count) {
            $nextText = $d.Paragraphs($i + 1).Range.Text
            if ($nextText.Length -gt 0 -and [int][char]$nextText[0] -eq 12) {
                $target = $p
            }
        }
    }
}
if ($target -eq $null) {
    $target = $fallback
}

# ------------------------------------------------------------------
# 3. Replace the (empty) paragraph content with the three runs and
#    drop the "Note" paragraph style (falls back to the document
#    default style, i.e. no explicit pStyle is written).
# ------------------------------------------------------------------
$target.Style = $d.Styles("Normal")

$para = $target.Range
$paraStart = $para.Start

$run1 = "This is an "
$run2 = "inline note, such as inside a table"
$apos = [char]0x2019
$run3 = ". It" + $apos + "s a character style."

$para.Text = $run1 + $run2 + $run3

$noteRange = $d.Range($paraStart + $run1.Length, $paraStart + $run1.Length + $run2.Length)
$noteRange.Style = $d.Styles("InlineNote")

Write-Host "Edit complete"
